$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data
$ws.Range("D2").Value = "30.651.55"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "1.877.48"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'238.97"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.4802"
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("D8").Value = "'0.2828"
$ws.Range("E8").Value = "  -2.40%  "
$ws.Range("D9").Value = "'0.06519"
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").Value = "1.953.77"
$ws.Range("E10").Value = "  +3.18%  "
$ws.Range("D11").Value = "'0.07474"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Value = "'16.51"
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("D13").Value = "'5.095"
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("D14").Value = "'88.05"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "'0.6633"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "30.596.03"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "'13.28"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007581"
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.164.79"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").Value = "'227.84"
$ws.Range("E21").Value = "  +3.15%  "
$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'5.277"
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("D24").Value = "'6.152"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").Value = "'168.33"
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("D26").Value = "'9.279"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").Value = "'18.54"
$ws.Range("E27").Value = "  -3.51%  "
$ws.Range("D28").Value = "'1.934"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "'1.405"
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("D30").Value = "'0.09704"
$ws.Range("E30").Value = "  +5.22%  "
$ws.Range("D31").Value = "'4.343"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").Value = "'4.006"
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("D33").Value = "'0.05073"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "'1.224"
$ws.Range("E34").Value = "  +5.89%  "
$ws.Range("D35").Value = "'0.7493"
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("D36").Value = "'2.715"
$ws.Range("D37").Value = "'0.01863"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "'2.633"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").Value = "'0.9135"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").Value = "'106.04"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "'0.4270"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").Value = "'5.768"
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "'7.348"
$ws.Range("E45").Value = "  -3.78%  "
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("D47").Value = "'64.37"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("D48").Value = "'1.473"
$ws.Range("E48").Value = "  -8.70%  "
$ws.Range("D49").Value = "'8.911"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").Value = "'33.72"
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("D51").Value = "'0.05658"
$ws.Range("E51").Value = "  -1.11%  "
